$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-19 down to 17-20.
$ws.Rows.Item(16).Insert()

# Copy the date number format from the row above (D15) to the new D16 cell
# so the new date value displays/stores consistently with the rest of column D.
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 16 with the new weekly record.
$ws.Cells.Item(16, 1).Value = 11
$ws.Cells.Item(16, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(16, 3).Value = "Bíobío"
$ws.Cells.Item(16, 4).Value = 45202
$ws.Cells.Item(16, 5).Value = 8
$ws.Cells.Item(16, 6).Value = 300000000
$ws.Cells.Item(16, 7).Value = "Espárragos"
$ws.Cells.Item(16, 8).Value = "Verde"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 300
$ws.Cells.Item(16, 11).Value = 1600
$ws.Cells.Item(16, 12).Value = 1600
$ws.Cells.Item(16, 13).Value = 1600
$ws.Cells.Item(16, 14).Value = "$/kilo"
$ws.Cells.Item(16, 15).Value = "Provincia de Linares"
$ws.Cells.Item(16, 16).Value = 1600
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = "Hortaliza"
